$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 90 (this shifts existing rows 90..185 down to 91..186,
# carrying their formatting, i.e. old row 90 data ends up at row 91, etc.)
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new weekly record.
$ws.Cells.Item(90, 1).Value = 7
$ws.Cells.Item(90, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(90, 3).Value = "Ñuble"
$ws.Cells.Item(90, 4).Value = 44587
$ws.Cells.Item(90, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(90, 5).Value = 16
$ws.Cells.Item(90, 6).Value = 100112032
$ws.Cells.Item(90, 7).Value = "Zapallo italiano"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 100
$ws.Cells.Item(90, 11).Value = 8500
$ws.Cells.Item(90, 12).Value = 9000
$ws.Cells.Item(90, 13).Value = 8750
$ws.Cells.Item(90, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(90, 15).Value = "Región del Maule"
$ws.Cells.Item(90, 16).Value = 146
$ws.Cells.Item(90, 17).Value = 60
$ws.Cells.Item(90, 18).Value = "Hortaliza"
